# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1), columns AD, AE, AF -> Wins, Losses, Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the rest of row 1 (bold/centered/bordered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record for every player row (2 through 46)
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 76
    $ws.Cells.Item($row, 31).Value = 85
    $ws.Cells.Item($row, 32).Value = 0
}
